{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph that separated it from the requisites line),\n// leaving \"LOM3206: Eletr\u00f4nica (Requisito)\" followed directly by the\n// single blank paragraph that precedes the final page-break paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two footer paragraphs by their text content.\nlet jupiterIdx = -1;\nlet copyrightIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (jupiterIdx === -1 && t.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIdx = i;\n  }\n  if (copyrightIdx === -1 && t.indexOf(\"Contact: luizeleno@usp.br\") !== -1) {\n    copyrightIdx = i;\n  }\n}\n\nif (jupiterIdx === -1 || copyrightIdx === -1) {\n  throw new Error(\"Could not locate footer paragraphs to remove.\");\n}\n\n// The blank paragraph immediately before the \"Ver no Jupiter\" paragraph\n// (separating it from \"LOM3206: Eletr\u00f4nica (Requisito)\") is removed too.\nlet blankIdx = -1;\nif (jupiterIdx - 1 >= 0 && items[jupiterIdx - 1].text === \"\") {\n  blankIdx = jupiterIdx - 1;\n}\n\n// Delete from the bottom up so earlier indices stay valid.\nconst toDelete = [copyrightIdx, jupiterIdx];\nif (blankIdx !== -1) {\n  toDelete.push(blankIdx);\n}\ntoDelete.sort((a, b) => b - a);\n\nfor (const idx of toDelete) {\n  items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n# (and the blank paragraph that separated it from the requisites line),\n# leaving \"LOM3206: Eletronica (Requisito)\" followed directly by the\n# single blank paragraph that precedes the final page-break paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Ver no Jupiter ...\" paragraph.\n$jupiter = $d.Content\n$jupiter.Find.Execute(\"Ver no Jupiter\") | Out-Null\n$jupiter.Expand(4)  # wdParagraph - grows the range to the whole paragraph (incl. mark)\n\n# Locate the \"(c) 2020 ... Contact: luizeleno@usp.br ...\" paragraph.\n$copyright = $d.Content\n$copyright.Find.Execute(\"Contact: luizeleno\") | Out-Null\n$copyright.Expand(4)  # wdParagraph\n\n# Extend one character to the left so the deletion also swallows the\n# blank paragraph mark that sits right before \"Ver no Jupiter ...\".\n$deleteRange = $d.Range($jupiter.Start - 1, $copyright.End)\n$deleteRange.Delete()\n"}
